$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "food" expense as the new row 2 (pushes old rows 2-4 down to 3-5) ---
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "food"
$ws.Range("B2").Value = 4000
# Copy the date cell formatting (numFmt) from the row below so the new date cell
# keeps the same date style used throughout column C.
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C2").Value = 45755.22928240741

# --- Insert a new "Buy Cloths" expense as the new row 4 (pushes old rows 4-5 down to 5-6) ---
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "Buy Cloths"
$ws.Range("B4").Value = 5000
$ws.Range("C5").Copy()
$ws.Range("C4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C4").Value = 45737.22928240741
